$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.6279234177358902
$ws.Range("F2").Value = 1.27351935105921

$ws.Range("C3").Value = 0.6279234177358902
$ws.Range("D3").Value = 0.003528628827193514
$ws.Range("E3").Value = 0.005619520991774325
$ws.Range("F3").Value = 1.27351935105921
$ws.Range("G3").Value = 0.01043734166128485
$ws.Range("H3").Value = 0.008195667896686389

$ws.Range("C4").Value = 0.6279234177358902
$ws.Range("D4").Value = 0.007946586843113446
$ws.Range("E4").Value = 0.0126553439777203
$ws.Range("F4").Value = 1.27351935105921
$ws.Range("G4").Value = 0.02353525321227257
$ws.Range("H4").Value = 0.01848048338857031

$ws.Range("C5").Value = 0.6279234177358902
$ws.Range("D5").Value = 0.01606891090935042
$ws.Range("E5").Value = 0.02559055842715701
$ws.Range("F5").Value = 1.27351935105921
$ws.Range("G5").Value = 0.04887858320107918
$ws.Range("H5").Value = 0.03838071495374291

$ws.Range("C6").Value = 0.6279234177358902
$ws.Range("D6").Value = 0.02754260422669819
$ws.Range("E6").Value = 0.04386299897208617
$ws.Range("F6").Value = 1.27351935105921
$ws.Range("G6").Value = 0.08242379185190507
$ws.Range("H6").Value = 0.06472127163466472

$ws.Range("C7").Value = 0.6279234177358902
$ws.Range("D7").Value = 0.04003080472359422
$ws.Range("E7").Value = 0.06375109383232384
$ws.Range("F7").Value = 1.27351935105921
$ws.Range("G7").Value = 0.1143468201407329
$ws.Range("H7").Value = 0.08978805076313018

$ws.Range("C8").Value = 0.6279234177358902
$ws.Range("D8").Value = 0.07342140686978549
$ws.Range("E8").Value = 0.1169273271166121
$ws.Range("F8").Value = 1.27351935105921
$ws.Range("G8").Value = 0.2044475054171753
$ws.Range("H8").Value = 0.1605374156640277

$ws.Range("C9").Value = 0.6279234177358902
$ws.Range("D9").Value = 0.1726146320864818
$ws.Range("E9").Value = 0.2748975865701587
$ws.Range("F9").Value = 1.27351935105921
$ws.Range("G9").Value = 0.4331075750530472
$ws.Range("H9").Value = 0.3400871566598682
